# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go count) figures on the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 673
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 8
$ws1.Range("F13").Value = 307
$ws1.Range("F14").Value = 416
$ws1.Range("F15").Value = 493
$ws1.Range("F16").Value = 128
$ws1.Range("F17").Value = 11200
$ws1.Range("F18").Value = 5341

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 92

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 673
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 92
$ws4.Range("F12").Value = 106
$ws4.Range("F13").Value = 8
$ws4.Range("F15").Value = 307
$ws4.Range("F16").Value = 416
$ws4.Range("F17").Value = 493
$ws4.Range("F18").Value = 128
$ws4.Range("F19").Value = 11200
$ws4.Range("F21").Value = 5341
